$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the image_width (column I) values to match the "old site" sizing.
$ws.Range("I3").Value = 100
$ws.Range("I4").Value = 50
$ws.Range("I8").Value = 33
$ws.Range("I13").Value = 50
$ws.Range("I16").Value = 50
$ws.Range("I16").Font.Bold = $true
$ws.Range("I20").Value = 50
$ws.Range("I33").Value = 50

# Match the saved selection state.
$ws.Range("I33").Select()
